$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Formula = "=B2/(12*1000)"
$ws.Range("B5").Formula = "=B3/(12*1000)"
$ws.Range("B6").Value = 0.080189449009009
$ws.Range("B7").Value = 43.773064159559
